# The "Förändrad" (changed) date in column C was bumped by one day
# (2023-09-20 -> 2023-09-21, i.e. Excel serial 45189 -> 45190) for every
# data row of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45190
}
